$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.0,1.0,0.9458723333333334,2.837617,0.1873686327665471,0.1873686327665471,3.0,1.0,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,41.32172242094778,371.8955017885301,0.08723827688277072,0.08723827688277072),
    @(3.0,1.0,0.9458723333333334,2.837617,0.1873686327665471,0.1873686327665471,3.0,1.0,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,1.139969750204222,10.259727751838,0.002406700178013004,0.002406700178013003),
    @(3.0,1.0,0.9458723333333334,2.837617,0.1873686327665471,0.1873686327665471,3.0,1.0,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,4.711602283026778,42.404420547241,0.009947118378584691,0.009947118378584691),
    @(3.0,1.0,0.9458723333333334,2.837617,0.1873686327665471,0.1873686327665471,3.0,1.0,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,41.57667757903612,374.1900982113251,0.08777653732717867,0.08777653732717866),
    @(3.0,1.0,1.964334,5.893002,0.3891165466060174,0.3891165466060174,3.0,1.0,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,85.81460883202,772.33147948818,0.1811715041694216,0.1811715041694216),
    @(3.0,1.0,1.964334,5.893002,0.3891165466060174,0.3891165466060174,3.0,1.0,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,2.367424503692,21.306820533228,0.004998098391160957,0.004998098391160957),
    @(3.0,1.0,1.964334,5.893002,0.3891165466060174,0.3891165466060174,3.0,1.0,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,9.784788319594,88.063094876346,0.02065761112202117,0.02065761112202117),
    @(3.0,1.0,1.964334,5.893002,0.3891165466060174,0.3891165466060174,3.0,1.0,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,86.34408524005,777.0967671604501,0.1822893329234137,0.1822893329234137),
    @(3.0,1.0,1.499502333333333,4.498507,0.2970376573303378,0.2970376573303378,3.0,1.0,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,65.50780375318111,589.57023377863,0.1382998478036614,0.1382998478036614),
    @(3.0,1.0,1.499502333333333,4.498507,0.2970376573303378,0.2970376573303378,3.0,1.0,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,1.807207209810889,16.264864888298,0.003815369585709678,0.003815369585709678),
    @(3.0,1.0,1.499502333333333,4.498507,0.2970376573303378,0.2970376573303378,3.0,1.0,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,7.469357510690111,67.224217596211,0.01576928163874542,0.01576928163874543),
    @(3.0,1.0,1.499502333333333,4.498507,0.2970376573303378,0.2970376573303378,3.0,1.0,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,65.91198710961945,593.207883986575,0.1391531583022214,0.1391531583022214),
    @(3.0,1.0,0.6384806666666667,1.915442,0.1264771632970977,0.1264771632970977,3.0,1.0,43.68636333333333,131.05909,0.4655970190670373,0.4655970190670373,27.89289838530889,251.03608546778,0.05888739021118358,0.05888739021118358),
    @(3.0,1.0,0.6384806666666667,1.915442,0.1264771632970977,0.1264771632970977,3.0,1.0,1.205204666666667,3.615614,0.01284473362738172,0.01284473362738172,0.7694998790431111,6.925498911388,0.001624565472498079,0.001624565472498079),
    @(3.0,1.0,0.6384806666666667,1.915442,0.1264771632970977,0.1264771632970977,3.0,1.0,4.981224333333333,14.943673,0.05308849315764798,0.05308849315764799,3.180415433162889,28.623738898466,0.006714482018296697,0.006714482018296697),
    @(3.0,1.0,0.6384806666666667,1.915442,0.1264771632970977,0.1264771632970977,3.0,1.0,43.95590833333333,131.867725,0.4684697541479331,0.4684697541479331,28.06499765660556,252.58497890945,0.05925072559511933,0.05925072559511934)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = $j + 5
        $ws.Cells.Item($row, $col).Value2 = $vals[$j]
    }
}
